$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume cells remain plain text (matches original inlineStr cells),
# so values like "7.00" or "5.80" do not get auto-converted to numbers and lose
# their trailing zeros / formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.633.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.762.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.85"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0684"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.018.76"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.19"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.758.38"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.655.06"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.607"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.09"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0769"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.55"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.02"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.75"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.06"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0508"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.48"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.77"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.379.53"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.651"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "77.59"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.42"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +13.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0142"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +18.18%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0499"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.29"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.80"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.917.39"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.49%  "
